$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header and value for PrefaultTime field in column G
$ws.Range("G1").Value = "PrefaultTime"

# Leading apostrophe forces the numeric-looking value to be stored as text
# (matching the existing quote-prefixed text style used by F2).
$ws.Range("G2").Value = "'2569"

# Move the active selection to the next empty cell, one column past the
# newly-added data (mirrors where Excel would leave the cursor after typing).
[void]$ws.Range("H2").Select()
